$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of quotes (2025-09-24) to the table
$ws.Range("A20").Value = 45924
$ws.Range("A20").NumberFormat = $ws.Range("A19").NumberFormat

$ws.Range("B20").Value = "20,9721"
$ws.Range("C20").Value = "15,1226"
$ws.Range("D20").Value = "14,9469"
$ws.Range("E20").Value = "14,9469"
